$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "sound for tire squeaking"
$ws.Range("C13").Select()
